$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.Value = '''59.320.53'
$c.Style = "Normal"

$c = $ws.Range('E2')
$c.Value = '''  +2.38%  '
$c.Style = "Normal"

$c = $ws.Range('D3')
$c.Value = '''3.167.56'
$c.Style = "Normal"

$c = $ws.Range('E3')
$c.Value = '''  +1.28%  '
$c.Style = "Normal"

$c = $ws.Range('E4')
$c.Value = '''  +0.01%  '
$c.Style = "Normal"

$c = $ws.Range('D5')
$c.Value = '''532.42'
$c.Style = "Normal"

$c = $ws.Range('E5')
$c.Value = '''  -0.25%  '
$c.Style = "Normal"

$c = $ws.Range('D6')
$c.Value = '''143.29'
$c.Style = "Normal"

$c = $ws.Range('E6')
$c.Value = '''  +3.04%  '
$c.Style = "Normal"

$c = $ws.Range('E7')
$c.Value = '''  +0.21%  '
$c.Style = "Normal"

$c = $ws.Range('D8')
$c.Value = '''0.515'
$c.Style = "Normal"

$c = $ws.Range('E8')
$c.Value = '''  +6.72%  '
$c.Style = "Normal"

$c = $ws.Range('D9')
$c.Value = '''7.25'
$c.Style = "Normal"

$c = $ws.Range('E9')
$c.Value = '''  -1.52%  '
$c.Style = "Normal"

$c = $ws.Range('D10')
$c.Value = '''0.112'
$c.Style = "Normal"

$c = $ws.Range('E10')
$c.Value = '''  +4.05%  '
$c.Style = "Normal"

$c = $ws.Range('D11')
$c.Value = '''0.430'
$c.Style = "Normal"

$c = $ws.Range('E11')
$c.Value = '''  +3.90%  '
$c.Style = "Normal"

$c = $ws.Range('D12')
$c.Value = '''3.713.37'
$c.Style = "Normal"

$c = $ws.Range('E12')
$c.Value = '''  +1.33%  '
$c.Style = "Normal"

$c = $ws.Range('D13')
$c.Value = '''0.139'
$c.Style = "Normal"

$c = $ws.Range('E13')
$c.Value = '''  +0.54%  '
$c.Style = "Normal"

$c = $ws.Range('D14')
$c.Value = '''25.83'
$c.Style = "Normal"

$c = $ws.Range('E14')
$c.Value = '''  -0.27%  '
$c.Style = "Normal"

$c = $ws.Range('D15')
$c.Value = '''0.0000172'
$c.Style = "Normal"

$c = $ws.Range('E15')
$c.Value = '''  +3.59%  '
$c.Style = "Normal"

$c = $ws.Range('D16')
$c.Value = '''59.344.77'
$c.Style = "Normal"

$c = $ws.Range('E16')
$c.Value = '''  +2.24%  '
$c.Style = "Normal"

$c = $ws.Range('D17')
$c.Value = '''3.165.81'
$c.Style = "Normal"

$c = $ws.Range('E17')
$c.Value = '''  +1.29%  '
$c.Style = "Normal"

$c = $ws.Range('D18')
$c.Value = '''6.17'
$c.Style = "Normal"

$c = $ws.Range('E18')
$c.Value = '''  +1.21%  '
$c.Style = "Normal"

$c = $ws.Range('D19')
$c.Value = '''12.92'
$c.Style = "Normal"

$c = $ws.Range('E19')
$c.Value = '''  +0.74%  '
$c.Style = "Normal"

$c = $ws.Range('D20')
$c.Value = '''8.10'
$c.Style = "Normal"

$c = $ws.Range('E20')
$c.Value = '''  -0.40%  '
$c.Style = "Normal"

$c = $ws.Range('D21')
$c.Value = '''373.75'
$c.Style = "Normal"

$c = $ws.Range('E21')
$c.Value = '''  -0.56%  '
$c.Style = "Normal"

$c = $ws.Range('D22')
$c.Value = '''0.998'
$c.Style = "Normal"

$c = $ws.Range('E22')
$c.Value = '''  -0.12%  '
$c.Style = "Normal"

$c = $ws.Range('D23')
$c.Value = '''0.526'
$c.Style = "Normal"

$c = $ws.Range('E23')
$c.Value = '''  +3.73%  '
$c.Style = "Normal"

$c = $ws.Range('D24')
$c.Value = '''69.81'
$c.Style = "Normal"

$c = $ws.Range('E24')
$c.Value = '''  +0.30%  '
$c.Style = "Normal"

$c = $ws.Range('E25')
$c.Value = '''  +0.63%  '
$c.Style = "Normal"

$c = $ws.Range('D26')
$c.Value = '''8.65'
$c.Style = "Normal"

$c = $ws.Range('E26')
$c.Value = '''  +15.77%  '
$c.Style = "Normal"

$c = $ws.Range('E27')
$c.Value = '''  -0.03%  '
$c.Style = "Normal"

$c = $ws.Range('D28')
$c.Value = '''0.0₃0884'
$c.Style = "Normal"

$c = $ws.Range('E28')
$c.Value = '''  +0.50%  '
$c.Style = "Normal"

$c = $ws.Range('E29')
$c.Value = '''  +1.26%  '
$c.Style = "Normal"

$c = $ws.Range('D30')
$c.Value = '''22.21'
$c.Style = "Normal"

$c = $ws.Range('E30')
$c.Value = '''  +3.17%  '
$c.Style = "Normal"

$c = $ws.Range('D31')
$c.Value = '''6.08'
$c.Style = "Normal"

$c = $ws.Range('E31')
$c.Value = '''  -1.81%  '
$c.Style = "Normal"

$c = $ws.Range('D32')
$c.Value = '''5.34'
$c.Style = "Normal"

$c = $ws.Range('E32')
$c.Value = '''  +3.71%  '
$c.Style = "Normal"

$c = $ws.Range('E33')
$c.Value = '''  -1.25%  '
$c.Style = "Normal"

$c = $ws.Range('D34')
$c.Value = '''6.38'
$c.Style = "Normal"

$c = $ws.Range('E34')
$c.Value = '''  +3.50%  '
$c.Style = "Normal"

$c = $ws.Range('D35')
$c.Value = '''156.36'
$c.Style = "Normal"

$c = $ws.Range('E35')
$c.Value = '''  -2.55%  '
$c.Style = "Normal"

$c = $ws.Range('D36')
$c.Value = '''1.34'
$c.Style = "Normal"

$c = $ws.Range('E36')
$c.Value = '''  +3.55%  '
$c.Style = "Normal"

$c = $ws.Range('D37')
$c.Value = '''0.0713'
$c.Style = "Normal"

$c = $ws.Range('E37')
$c.Value = '''  +5.99%  '
$c.Style = "Normal"

$c = $ws.Range('D38')
$c.Value = '''25.42'
$c.Style = "Normal"

$c = $ws.Range('E38')
$c.Value = '''  -0.88%  '
$c.Style = "Normal"

$c = $ws.Range('D39')
$c.Value = '''2.717.18'
$c.Style = "Normal"

$c = $ws.Range('E39')
$c.Value = '''  +6.78%  '
$c.Style = "Normal"

$c = $ws.Range('D40')
$c.Value = '''1.66'
$c.Style = "Normal"

$c = $ws.Range('E40')
$c.Value = '''  +0.28%  '
$c.Style = "Normal"

$c = $ws.Range('D41')
$c.Value = '''4.25'
$c.Style = "Normal"

$c = $ws.Range('E41')
$c.Value = '''  +3.40%  '
$c.Style = "Normal"

$c = $ws.Range('B42')
$c.Value = '''OKB'
$c.Style = "Normal"

$c = $ws.Range('C42')
$c.Value = '''https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c.Style = "Normal"

$c = $ws.Range('D42')
$c.Value = '''39.34'
$c.Style = "Normal"

$c = $ws.Range('E42')
$c.Value = '''  +3.88%  '
$c.Style = "Normal"

$c = $ws.Range('B43')
$c.Value = '''Mantle'
$c.Style = "Normal"

$c = $ws.Range('C43')
$c.Value = '''https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$c.Style = "Normal"

$c = $ws.Range('D43')
$c.Value = '''0.721'
$c.Style = "Normal"

$c = $ws.Range('E43')
$c.Value = '''  +2.88%  '
$c.Style = "Normal"

$c = $ws.Range('B44')
$c.Value = '''VeChain'
$c.Style = "Normal"

$c = $ws.Range('C44')
$c.Value = '''https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c.Style = "Normal"

$c = $ws.Range('D44')
$c.Value = '''0.0290'
$c.Style = "Normal"

$c = $ws.Range('E44')
$c.Value = '''  +7.21%  '
$c.Style = "Normal"

$c = $ws.Range('D45')
$c.Value = '''3.211.94'
$c.Style = "Normal"

$c = $ws.Range('E45')
$c.Value = '''  +1.39%  '
$c.Style = "Normal"

$c = $ws.Range('D46')
$c.Value = '''0.987'
$c.Style = "Normal"

$c = $ws.Range('E46')
$c.Value = '''  +0.21%  '
$c.Style = "Normal"

$c = $ws.Range('D47')
$c.Value = '''6.16'
$c.Style = "Normal"

$c = $ws.Range('E47')
$c.Value = '''  -0.07%  '
$c.Style = "Normal"

$c = $ws.Range('D48')
$c.Value = '''0.0998'
$c.Style = "Normal"

$c = $ws.Range('E48')
$c.Value = '''  +7.87%  '
$c.Style = "Normal"

$c = $ws.Range('D49')
$c.Value = '''20.40'
$c.Style = "Normal"

$c = $ws.Range('E49')
$c.Value = '''  +2.84%  '
$c.Style = "Normal"

$c = $ws.Range('E50')
$c.Value = '''  -0.02%  '
$c.Style = "Normal"

$c = $ws.Range('D51')
$c.Value = '''0.764'
$c.Style = "Normal"

$c = $ws.Range('E51')
$c.Value = '''  +1.81%  '
$c.Style = "Normal"

